$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: mirrors row 1's "columnXY" cells (F1,G1,H1,I1), producing the
# "clear the cell if its text is the string 0" statement used by btnDown.
$ws.Range("F21").Formula = '="if("&F1&".getText().equals("&CHAR(34)&"0"&CHAR(34)&")) "&F1&".setText(null);"'
$ws.Range("G21:I21").Formula = '="if("&G1&".getText().equals("&CHAR(34)&"0"&CHAR(34)&")) "&G1&".setText(null);"'

# Rows 22-24: same pattern anchored on row 2's column names (F2,G2,H2,I2),
# filled down as one shared formula block across F22:I24.
$ws.Range("F22:I24").Formula = '="if("&F2&".getText().equals("&CHAR(34)&"0"&CHAR(34)&")) "&F2&".setText(null);"'

# Restore the selection to the newly added block, as the author left it.
$ws.Range("F21:I24").Select() | Out-Null
